# "all members sprint task assgined for sprint1"
# Adds the Sprint1 task rows for US03 (Birth before death) and
# US06 (Divorce before death), owned by "myl", below the existing
# US02 task block on the Sprint1 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# --- Row 28: US03 story summary row -----------------------------------
$ws.Range("A28").Value = "US03"
$ws.Range("B28").Value = "Birth before death"
$ws.Range("C28").Value = "myl"
$ws.Range("D28").Value = "coding"
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 60

# --- Row 30-32: US03 task breakdown -------------------------------------
$ws.Range("A30").Value = "T03.01"
$ws.Range("B30").Value = "Store birth date"
$ws.Range("C30").Value = "myl"
$ws.Range("D30").Value = "Done"

$ws.Range("A31").Value = "T03.02"
$ws.Range("B31").Value = "Store death date"
$ws.Range("C31").Value = "myl"
$ws.Range("D31").Value = "Done"

$ws.Range("A32").Value = "T03.02"
$ws.Range("B32").Value = "Compare to birth date and death date"
$ws.Range("C32").Value = "myl"

# --- Row 34: US06 story summary row -------------------------------------
$ws.Range("A34").Value = "US06"
$ws.Range("B34").Value = "Divorce before death"
$ws.Range("C34").Value = "myl"
$ws.Range("D34").Value = "Coding"
$ws.Range("E34").Value = 150
$ws.Range("F34").Value = 80

# --- Row 36-38: US06 task breakdown -------------------------------------
$ws.Range("A36").Value = "T06.01"
$ws.Range("B36").Value = "Store divorce date"
$ws.Range("C36").Value = "myl"

$ws.Range("A37").Value = "T06.02"
$ws.Range("B37").Value = "Store death date"
$ws.Range("C37").Value = "myl"

$ws.Range("A38").Value = "T06.03"
$ws.Range("B38").Value = "Cpmpare to divorce date and death date"
$ws.Range("C38").Value = "myl"

# --- Formatting: the whole new block (rows 28-38, cols A-G) picks up an
# explicit black font colour (new font/style in styles.xml), matching the
# newly typed-in block in the source edit.
$ws.Range("A28:G32").Font.Color = 0
$ws.Range("A33:G38").Font.Color = 0

# --- View state: scrolled down a bit, with C41 selected afterwards -----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C41").Select()
